$wb = $excel.ActiveWorkbook

# ---- Sheet1 (Schedule) updates ----
$ws1 = $wb.Worksheets.Item("Schedule")

# Row 2 (existing row) - update values
$ws1.Cells.Item(2,2).Value = 46040.83333333334
$ws1.Cells.Item(2,3).Value = 13.5
$ws1.Cells.Item(2,4).Value = 51.02999999999999
$ws1.Cells.Item(2,5).Value = 60.30140025000001
$ws1.Cells.Item(2,6).Value = 1.181685288065844

# New rows 3 and 4 - set values and copy date style from row 2 for columns A/B
# Row 3
$ws1.Cells.Item(3,1).Value = 46040.89583333334
$ws1.Cells.Item(3,2).Value = 46041.125
$ws1.Cells.Item(3,3).Value = 5.5
$ws1.Cells.Item(3,4).Value = 20.79
$ws1.Cells.Item(3,5).Value = 450.2433195
$ws1.Cells.Item(3,6).Value = 21.65672532467532
$ws1.Cells.Item(3,1).NumberFormat = $ws1.Cells.Item(2,1).NumberFormat
$ws1.Cells.Item(3,2).NumberFormat = $ws1.Cells.Item(2,2).NumberFormat

# Row 4
$ws1.Cells.Item(4,1).Value = 46041.29166666666
$ws1.Cells.Item(4,2).Value = 46041.66666666666
$ws1.Cells.Item(4,3).Value = 9
$ws1.Cells.Item(4,4).Value = 34.02
$ws1.Cells.Item(4,5).Value = -43.64514375000001
$ws1.Cells.Item(4,6).Value = -1.282926036155203
$ws1.Cells.Item(4,1).NumberFormat = $ws1.Cells.Item(2,1).NumberFormat
$ws1.Cells.Item(4,2).NumberFormat = $ws1.Cells.Item(2,2).NumberFormat

# ---- Sheet2 (Detailed) updates ----
$ws2 = $wb.Worksheets.Item("Detailed")

# --- Modify existing rows 25-49 ---
$ws2.Cells.Item(25,2).Value = -5.33031
$ws2.Cells.Item(26,2).Value = 0
$ws2.Cells.Item(27,2).Value = -5.01
$ws2.Cells.Item(27,3).Value = "historical"
$ws2.Cells.Item(28,2).Value = -5.51
$ws2.Cells.Item(28,3).Value = "historical"
$ws2.Cells.Item(29,2).Value = -5.58988
$ws2.Cells.Item(29,3).Value = "historical"
$ws2.Cells.Item(30,3).Value = "historical"
$ws2.Cells.Item(31,2).Value = -9.99
$ws2.Cells.Item(31,3).Value = "historical"
$ws2.Cells.Item(32,2).Value = -8.452769999999999
$ws2.Cells.Item(33,2).Value = -7.62043
$ws2.Cells.Item(34,2).Value = -6.32
$ws2.Cells.Item(35,2).Value = -6.34762
$ws2.Cells.Item(36,2).Value = -5.95309
$ws2.Cells.Item(37,2).Value = 4.94039
$ws2.Cells.Item(38,2).Value = 6.91421
$ws2.Cells.Item(39,2).Value = 12.15556
$ws2.Cells.Item(40,2).Value = 27.49695
$ws2.Cells.Item(41,2).Value = 36.0601
$ws2.Cells.Item(42,2).Value = 45.40186
$ws2.Cells.Item(42,5).Value = "OFF"
$ws2.Cells.Item(43,2).Value = 56.98
$ws2.Cells.Item(43,5).Value = "OFF"
$ws2.Cells.Item(44,2).Value = 36.2
$ws2.Cells.Item(44,5).Value = "OFF"
$ws2.Cells.Item(46,2).Value = 46.4787
$ws2.Cells.Item(47,5).Value = "ON"
$ws2.Cells.Item(48,5).Value = "ON"
$ws2.Cells.Item(49,5).Value = "ON"

# --- Add new rows 50-97 ---
# Data array: DateTime, Price, Type, Date, Pump_Status
$newRows = @(
    @(50, 46041, 36.06, "forecast", 46041, "ON"),
    @(51, 46041.02083333334, 36.06, "forecast", 46041, "ON"),
    @(52, 46041.04166666666, 57.06003, "forecast", 46041, "ON"),
    @(53, 46041.0625, 56.98, "forecast", 46041, "ON"),
    @(54, 46041.08333333334, 36.06, "forecast", 46041, "ON"),
    @(55, 46041.10416666666, 48.84899, "forecast", 46041, "ON"),
    @(56, 46041.125, 49.74012, "forecast", 46041, "OFF"),
    @(57, 46041.14583333334, 56.98, "forecast", 46041, "OFF"),
    @(58, 46041.16666666666, 56.98, "forecast", 46041, "OFF"),
    @(59, 46041.1875, 57.06003, "forecast", 46041, "OFF"),
    @(60, 46041.20833333334, 58.32151, "forecast", 46041, "OFF"),
    @(61, 46041.22916666666, 69.20653, "forecast", 46041, "OFF"),
    @(62, 46041.25, 72.94628, "forecast", 46041, "OFF"),
    @(63, 46041.27083333334, 57.06003, "forecast", 46041, "OFF"),
    @(64, 46041.29166666666, 30.36901, "forecast", 46041, "ON"),
    @(65, 46041.3125, 6.75696, "forecast", 46041, "ON"),
    @(66, 46041.33333333334, 0.7, "forecast", 46041, "ON"),
    @(67, 46041.35416666666, 0.64751, "forecast", 46041, "ON"),
    @(68, 46041.375, 0, "forecast", 46041, "ON"),
    @(69, 46041.39583333334, -6.05409, "forecast", 46041, "ON"),
    @(70, 46041.41666666666, -7.79943, "forecast", 46041, "ON"),
    @(71, 46041.4375, -6.97876, "forecast", 46041, "ON"),
    @(72, 46041.45833333334, -7.69821, "forecast", 46041, "ON"),
    @(73, 46041.47916666666, -6.90017, "forecast", 46041, "ON"),
    @(74, 46041.5, -7.02264, "forecast", 46041, "ON"),
    @(75, 46041.52083333334, -6.94073, "forecast", 46041, "ON"),
    @(76, 46041.54166666666, -6.08034, "forecast", 46041, "ON"),
    @(77, 46041.5625, -5.95857, "forecast", 46041, "ON"),
    @(78, 46041.58333333334, -5.01, "forecast", 46041, "ON"),
    @(79, 46041.60416666666, -5.51, "forecast", 46041, "ON"),
    @(80, 46041.625, -5.77494, "forecast", 46041, "ON"),
    @(81, 46041.64583333334, -5.50985, "forecast", 46041, "ON"),
    @(82, 46041.66666666666, -2.52431, "forecast", 46041, "OFF"),
    @(83, 46041.6875, -5.13343, "forecast", 46041, "OFF"),
    @(84, 46041.70833333334, -2.63766, "forecast", 46041, "OFF"),
    @(85, 46041.72916666666, 2.2195, "forecast", 46041, "OFF"),
    @(86, 46041.75, 20.24437, "forecast", 46041, "OFF"),
    @(87, 46041.77083333334, 47.97146, "forecast", 46041, "OFF"),
    @(88, 46041.79166666666, 59.14369, "forecast", 46041, "OFF"),
    @(89, 46041.8125, 77.94, "forecast", 46041, "OFF"),
    @(90, 46041.83333333334, 73.20007, "forecast", 46041, "OFF"),
    @(91, 46041.85416666666, 66.89706, "forecast", 46041, "OFF"),
    @(92, 46041.875, 67.11141000000001, "forecast", 46041, "OFF"),
    @(93, 46041.89583333334, 62.40577, "forecast", 46041, "OFF"),
    @(94, 46041.91666666666, 57.54921, "forecast", 46041, "OFF"),
    @(95, 46041.9375, 57.3, "forecast", 46041, "OFF"),
    @(96, 46041.95833333334, 57.06007, "forecast", 46041, "OFF"),
    @(97, 46041.97916666666, 57.06003, "forecast", 46041, "OFF")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws2.Cells.Item($r,1).Value = $row[1]
    $ws2.Cells.Item($r,2).Value = $row[2]
    $ws2.Cells.Item($r,3).Value = $row[3]
    $ws2.Cells.Item($r,4).Value = $row[4]
    $ws2.Cells.Item($r,5).Value = $row[5]
    $ws2.Cells.Item($r,1).NumberFormat = $ws2.Cells.Item(2,1).NumberFormat
    $ws2.Cells.Item($r,4).NumberFormat = $ws2.Cells.Item(2,4).NumberFormat
}